# Finished grading weeks 2 and 3
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Week 3 grading must be entered before Week 2 so that the two brand-new
# shared strings land in the same order as the target workbook:
#   "Failed to run"  (Program_03_6 / D12) -> first new shared string
#   "Does not run"   (Tutorial_03_2 / D3) -> second new shared string
# ---------------------------------------------------------------------------
$week3 = $wb.Worksheets.Item("Week 3")

# Row 12 - Program_03_6 failed to run
$week3.Range("B12").Value = 0
$week3.Range("D12").Value = "Failed to run"

# Row 3 - Tutorial_03_2 does not run
$week3.Range("B3").Value = 0
$week3.Range("D3").Value = "Does not run"

# Remaining rows on Week 3 all graded full marks, no notes
$week3.Range("B2").Value = 1
$week3.Range("D2").Style = "Normal"

$week3.Range("B4").Value = 1
$week3.Range("D4").Style = "Normal"

$week3.Range("B5").Value = 1
$week3.Range("D5").Style = "Normal"

$week3.Range("B6").Value = 1
$week3.Range("D6").Style = "Normal"

$week3.Range("B7").Value = 1
$week3.Range("D7").Style = "Normal"

$week3.Range("B8").Value = 1
$week3.Range("D8").Style = "Normal"

$week3.Range("B9").Value = 1
$week3.Range("D9").Style = "Normal"

$week3.Range("B10").Value = 1
$week3.Range("D10").Style = "Normal"

$week3.Range("B11").Value = 1
$week3.Range("D11").Style = "Normal"

# Column widths settled back to their content-driven auto size
$week3.Range("B1").EntireColumn.ColumnWidth = 8.333333333333332
$week3.Range("D1").EntireColumn.ColumnWidth = 12.333333333333332

# ---------------------------------------------------------------------------
# Week 2 - all rows graded full marks, no notes
# ---------------------------------------------------------------------------
$week2 = $wb.Worksheets.Item("Week 2")

$week2.Range("B2").Value = 1
$week2.Range("D2").Style = "Normal"

$week2.Range("B3").Value = 1
$week2.Range("D3").Style = "Normal"

$week2.Range("B4").Value = 1
$week2.Range("D4").Style = "Normal"

$week2.Range("B5").Value = 1
$week2.Range("D5").Style = "Normal"

$week2.Range("B6").Value = 1
$week2.Range("D6").Style = "Normal"

$week2.Range("B7").Value = 1
$week2.Range("D7").Style = "Normal"

# Column widths settled back to their content-driven auto size
$week2.Range("B1").EntireColumn.ColumnWidth = 8.333333333333332
$week2.Range("D1").EntireColumn.ColumnWidth = 5.833333333333334
